$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-506). Every one of those cells moves from 45180 to 45181
# (i.e. the "changed" date advances by one day for the whole table).
$ws.Range("C2:C506").Value = 45181
